# Apply the "broadband dept" data update described by the commit.
$wb = $excel.ActiveWorkbook

$deptA = "Broadband, Communications and the Digital Economy portfolio"
$sourceUrl = "http://www.dbcde.gov.au/__data/assets/pdf_file/0004/163417/COMPLETE_BCDE_PBS.pdf"
$sourceTitle = "Broadband, Communications and the Digital Economy Portfolio Budget Statement"

# Column B (agency) / C (program) / F (13-14 estimate) / G (prior year actual)
$rows = @(
    @(754, "Department of Broadband, Communications and the Digital Economy", "Broadband and Communications Infrastructure", 70710, 48519),
    @(755, "Department of Broadband, Communications and the Digital Economy", "Digital Economy and Postal Services", 98301, 95983),
    @(756, "Department of Broadband, Communications and the Digital Economy", "Broadcasting and Digital Television", 268066, 210896),
    @(757, "Australian Broadcasting Corporation", "ABC Radio", 337691, 348571),
    @(758, "Australian Broadcasting Corporation", "ABC Television", 627150, 647348),
    @(759, "Australian Broadcasting Corporation", "Online", 25755, 26587),
    @(760, "Australian Broadcasting Corporation", "ABC Analog Transmission", 88669, 80403),
    @(761, "Australian Broadcasting Corporation", "Access to digital TV services", 100673, 109103),
    @(762, "Australian Broadcasting Corporation", "Access to digital radio services", 3658, 3705),
    @(763, "Australian Communications and Media Authority", "Communications regulation, planning and licensing", 73430, 68150),
    @(764, "Australian Communications and Media Authority", "Consumer safeguards, education and information", 86052, 25487),
    @(765, "Special Broadcasting Service Corporation", "Television", 214488, 250057),
    @(766, "Special Broadcasting Service Corporation", "Radio", 37518, 40474),
    @(767, "Special Broadcasting Service Corporation", "Analog Transmission and Distribution", 11017, 3420),
    @(768, "Special Broadcasting Service Corporation", "Digital TV Transmission and Distribution", 69069, 81440),
    @(769, "Special Broadcasting Service Corporation", "Digital Radio Transmission and Distribution", 2030, 2064),
    @(770, "Telecommunications Universal Service Management Agency", "", 313565, 346457)
)

# ---- "Raw data" sheet (columns A-L) ----
$raw = $wb.Worksheets.Item("Raw data")
foreach ($r in $rows) {
    $rowNum = $r[0]
    $raw.Cells.Item($rowNum, 1).Value = $deptA
    $raw.Cells.Item($rowNum, 2).Value = $r[1]
    if ($r[2] -ne "") {
        $raw.Cells.Item($rowNum, 3).Value = $r[2]
    }
    $raw.Cells.Item($rowNum, 6).Value = $r[3]
    $raw.Cells.Item($rowNum, 7).Value = $r[4]
    $raw.Cells.Item($rowNum, 11).Value = $sourceUrl
    $raw.Cells.Item($rowNum, 12).Value = $sourceTitle
}

# ---- "output sheet" (columns A-J, same data minus column H) ----
$out = $wb.Worksheets.Item("output sheet")
foreach ($r in $rows) {
    $rowNum = $r[0]
    $out.Cells.Item($rowNum, 1).Value = $deptA
    $out.Cells.Item($rowNum, 2).Value = $r[1]
    if ($r[2] -ne "") {
        $out.Cells.Item($rowNum, 3).Value = $r[2]
    }
    $out.Cells.Item($rowNum, 6).Value = $r[3]
    $out.Cells.Item($rowNum, 7).Value = $r[4]
    $out.Cells.Item($rowNum, 9).Value = $sourceUrl
    $out.Cells.Item($rowNum, 10).Value = $sourceTitle
}

# Restore the workbook's active tab/selection state on "output sheet".
$out.Activate()
